$wb = $excel.ActiveWorkbook

# Rename the "column_specs" sheet to "column_names" (column specs are now internal)
$ws = $wb.Worksheets.Item("column_specs")
$ws.Name = "column_names"

# Remove the now-unused "col_type" column from the renamed sheet's table
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item("col_type").Delete()

# Make the renamed sheet the active tab
$ws.Activate()
